$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC002_Output")

$ws.Range("A8").Value = "Dr. Chandrashekara Aithal"
$ws.Range("A9").Value = "Dr. J Prasad"
$ws.Range("A10").Value = "Dr. Roshan Kumar.B"
$ws.Range("A11").Value = "Dr. Revanth BN"
